$wb = $excel.ActiveWorkbook

# --- Fix 1: correct the LANGUAGES value on the "General" sheet ---
$general = $wb.Worksheets.Item("General")
$general.Range("B5").Value = 'en","da","kl'

# --- Fix 2: split/quote fix for code columns on the "Data" sheet ---
# The A and C columns held the long code-label (e.g. "Total", "Men",
# "Women", "Greenland", "Outside Greenland") instead of the short code
# (e.g. "T", "M", "K", "N", "S"). Remap them.
$data = $wb.Worksheets.Item("Data")

$map = @{
    "Total"              = "T"
    "Men"                = "M"
    "Women"              = "K"
    "Greenland"          = "N"
    "Outside Greenland"  = "S"
}

$lastRow = $data.Cells.Item($data.Rows.Count, 1).End(-4162).Row  # xlUp

for ($r = 2; $r -le $lastRow; $r++) {
    foreach ($col in 1, 3) {
        $cell = $data.Cells.Item($r, $col)
        $val = $cell.Value2
        if ($map.ContainsKey($val)) {
            $cell.Value = $map[$val]
        }
    }
}
